# Correction in SA algorithm and 746 logs
# Update Fitness (column C) values on Sheet1 for run_29 log data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-25 (generations 0-23): fitness corrected to 7982
$ws.Range("C2:C25").Value = 7982

# Rows 26-37 (generations 24-35): fitness corrected to 7682
$ws.Range("C26:C37").Value = 7682

# Rows 38-77 (generations 36-75): fitness corrected to 7657
$ws.Range("C38:C77").Value = 7657

# Rows 176-252 (generations 174-250): fitness corrected to 7573
$ws.Range("C176:C252").Value = 7573
